$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the STATUS / PASSWORD columns (E, F) — template now imports from
#     an Excel sheet instead of a CSV, so accounts are provisioned directly
#     from FIRST NAME / LAST NAME / EMAIL / GENDER. ---
$ws.Range("E1:F1").EntireColumn.Delete()

# --- Make the bold header (20% - Accent2 style used by the header row). ---
$wb.Styles.Item("20% - Accent2").Font.Bold = $true

# --- Sample data rows (two example teacher/student records) ---
$ws.Range("A2").Value = "lesala"
$ws.Range("B2").Value = "lesala"
$ws.Range("C2").Value = "lesala@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:lesala@gmail.com")
$ws.Range("D2").Value = "Male"

$ws.Range("A3").Value = "likobo"
$ws.Range("B3").Value = "likobo"
$ws.Range("C3").Value = "likobo@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:likobo@gmail.com")
$ws.Range("D3").Value = "Female"

# --- Drop the old STATUS (Active/Inactive) validation list; keep GENDER's. ---
$ws.Range("E1:E1048576").Validation.Delete()

# --- Update the active selection to reflect where the user ended up. ---
$null = $ws.Range("D5").Select()
